$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the 9 new test-data rows (A22:G30), following the same pattern
# as the existing rows.
$data = @(
    @(10002, 110021),
    @(10003, 110022),
    @(10004, 110023),
    @(10005, 110024),
    @(10006, 110025),
    @(10007, 110026),
    @(10008, 110027),
    @(10009, 110028),
    @(10010, 110029)
)

$row = 22
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $ws.Cells.Item($row, 3).Value = "eng"
    $ws.Cells.Item($row, 4).Value = $true
    $ws.Cells.Item($row, 5).Value = "superadmin"
    $ws.Cells.Item($row, 6).Value = "now()"
    $ws.Cells.Item($row, 7).Value = "now()"
    $row++
}

# Select the empty rows below the data (mirrors the author selecting
# everything from row 31 down before saving).
[void]$ws.Rows("31:1048576").Select()

# Explicitly set the page orientation to portrait.
$ws.PageSetup.Orientation = 1
